$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 08:37"

# Swap "Groenlandia" / "Islas Malvinas" text values (same underlying numbers)
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# Row 76 - El Salvador
$ws.Range("D76").Value = 7324
$ws.Range("E76").Value = 5674
$ws.Range("G76").Value = 7
$ws.Range("H76").Value = 379

# Row 101 - Hungria
$ws.Range("B101").Value = 4398
$ws.Range("C101").Value = 18
$ws.Range("D101").Value = 3312
$ws.Range("E101").Value = 490

# Row 144 - Georgia
$ws.Range("B144").Value = 1104
$ws.Range("C144").Value = 19
$ws.Range("D144").Value = 912
$ws.Range("E144").Value = 176

# Row 161 - Taiwan
$ws.Range("B161").Value = 458
$ws.Range("C161").Value = 3
$ws.Range("E161").Value = 11
